# Weekly fruit/vegetable update:
# Insert a new price record (row 120) for "Papa" / "Asterix" / "1a (cosecha)"
# in "Región del Maule", pushing the existing rows 120-139 down to 121-140.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 120; this shifts rows 120:139 down to 121:140,
# matching the growth of the sheet's dimension from A1:R139 to A1:R140.
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(120, 1).Value = 1
$ws.Cells.Item(120, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(120, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(120, 4).Value = 44951
$ws.Cells.Item(120, 5).Value = 15
$ws.Cells.Item(120, 6).Value = 100114001
$ws.Cells.Item(120, 7).Value = "Papa"
$ws.Cells.Item(120, 8).Value = "Asterix"
$ws.Cells.Item(120, 9).Value = "1a (cosecha)"
$ws.Cells.Item(120, 10).Value = 1000
$ws.Cells.Item(120, 11).Value = 14000
$ws.Cells.Item(120, 12).Value = 15000
$ws.Cells.Item(120, 13).Value = 14500
$ws.Cells.Item(120, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(120, 15).Value = "Región del Maule"
$ws.Cells.Item(120, 16).Value = 580
$ws.Cells.Item(120, 17).Value = 25
$ws.Cells.Item(120, 18).Value = "Hortaliza"
